$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that should move forward by
# one day (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data
# row (rows 2 through 321).
$ws.Range("C2:C321").Value = 45182
